{"js": "// Fill in the \"Exemple\" column for the `datetime` row of the `position`\n// table, and the `lat` / `lon` rows of the `coord` table \u2014 these cells\n// were previously empty.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < tables.items.length; i++) {\n  tables.items[i].load(\"values\");\n}\nawait context.sync();\n\n// Pairs of [rowLabel text to locate (col 0), example value to write (col 5)]\nconst targets = [\n  { row: \"datetime\", parentRowLabel: \"Date/heure de la position\", value: \"2022-09-27T08:23:34+02:00\" },\n  { row: \"lat\", parentRowLabel: \"Latitude\", value: \"48.866667\" },\n  { row: \"lon\", parentRowLabel: \"Longitude\", value: \"2.333333\" },\n];\n\nfor (const target of targets) {\n  for (let i = 0; i < tables.items.length; i++) {\n    const table = tables.items[i];\n    const rows = table.values;\n    for (let r = 1; r < rows.length; r++) {\n      const row = rows[r];\n      if (row[0] === target.row && row[1] === target.parentRowLabel && row[5] === \"\") {\n        const cell = table.getCell(r, 5);\n        cell.body.insertText(target.value, Word.InsertLocation.replace);\n      }\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Fill in the previously-empty \"Exemple\" column cells:\n#   - the `datetime` row of the `position` table\n#   - the `lat` and `lon` rows of the `coord` table\n\n$d = $word.ActiveDocument\n\n$targets = @(\n    @{ Row = \"datetime\"; Label = \"Date/heure de la position\"; Value = \"2022-09-27T08:23:34+02:00\" },\n    @{ Row = \"lat\";       Label = \"Latitude\";                  Value = \"48.866667\" },\n    @{ Row = \"lon\";       Label = \"Longitude\";                 Value = \"2.333333\" }\n)\n\nfor ($i = 1; $i -le $d.Tables.Count; $i++) {\n    $t = $d.Tables.Item($i)\n    for ($r = 2; $r -le $t.Rows.Count; $r++) {\n        $rowLabel = $t.Cell($r, 1).Range.Text.TrimEnd([char]7, [char]13)\n        $colLabel = $t.Cell($r, 2).Range.Text.TrimEnd([char]7, [char]13)\n        foreach ($target in $targets) {\n            if ($rowLabel -eq $target.Row -and $colLabel -eq $target.Label) {\n                $exCell = $t.Cell($r, 6)\n                $exText = $exCell.Range.Text.TrimEnd([char]7, [char]13)\n                if ([string]::IsNullOrEmpty($exText)) {\n                    $exCell.Range.Text = $target.Value\n                }\n            }\n        }\n    }\n}\n"}
